$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2488215.8
$ws.Cells.Item(86, 9).Value = 2500
$ws.Cells.Item(86, 10).Value = 3040597
$ws.Cells.Item(86, 11).Value = 2500
$ws.Cells.Item(86, 12).Value = 3040597
$ws.Cells.Item(86, 13).Value = -1377
$ws.Cells.Item(86, 14).Value = -3042843
$ws.Cells.Item(89, 8).Value = 2488215.8
$ws.Cells.Item(89, 9).Value = 2500
$ws.Cells.Item(89, 10).Value = 3040597
$ws.Cells.Item(89, 11).Value = 12500
$ws.Cells.Item(89, 12).Value = 15202985
$ws.Cells.Item(89, 13).Value = -6884
$ws.Cells.Item(89, 14).Value = -15214217
$ws.Cells.Item(137, 8).Value = 23453718
$ws.Cells.Item(137, 10).Value = 62549800
$ws.Cells.Item(137, 12).Value = 187649400
$ws.Cells.Item(137, 14).Value = -187654500

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7538098
$ws.Cells.Item(32, 9).Value = 1901370.6
$ws.Cells.Item(32, 11).Value = 1901370.6
$ws.Cells.Item(32, 13).Value = -1901083.6
$ws.Cells.Item(61, 8).Value = 5583762
$ws.Cells.Item(61, 9).Value = 2605199.5
$ws.Cells.Item(61, 10).Value = 29412264
$ws.Cells.Item(61, 11).Value = 2605199.5
$ws.Cells.Item(61, 12).Value = 29412264
$ws.Cells.Item(61, 13).Value = -2604987.5
$ws.Cells.Item(61, 14).Value = -29412688
$ws.Cells.Item(74, 8).Value = 45519740
$ws.Cells.Item(74, 9).Value = 39011748
$ws.Cells.Item(74, 10).Value = 66670704
$ws.Cells.Item(74, 11).Value = 39011748
$ws.Cells.Item(74, 12).Value = 66670704
$ws.Cells.Item(74, 13).Value = -39010874
$ws.Cells.Item(74, 14).Value = -66672452
$ws.Cells.Item(77, 8).Value = 45519740
$ws.Cells.Item(77, 9).Value = 39011748
$ws.Cells.Item(77, 10).Value = 66670704
$ws.Cells.Item(77, 11).Value = 195058740
$ws.Cells.Item(77, 12).Value = 333353520
$ws.Cells.Item(77, 13).Value = -195054372
$ws.Cells.Item(77, 14).Value = -333362256
$ws.Cells.Item(132, 8).Value = 19451052
$ws.Cells.Item(132, 9).Value = 22229966
$ws.Cells.Item(132, 10).Value = 11114311
$ws.Cells.Item(132, 11).Value = 66689898
$ws.Cells.Item(132, 12).Value = 33342933
$ws.Cells.Item(132, 13).Value = -66687368
$ws.Cells.Item(132, 14).Value = -33347993
$ws.Cells.Item(136, 8).Value = 5583762
$ws.Cells.Item(136, 9).Value = 2605199.5
$ws.Cells.Item(136, 10).Value = 29412264
$ws.Cells.Item(136, 11).Value = 7815598.5
$ws.Cells.Item(136, 12).Value = 88236792
$ws.Cells.Item(136, 13).Value = -7813048.5
$ws.Cells.Item(136, 14).Value = -88241892

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1896.43
$ws.Cells.Item(86, 9).Value = 1898.3265
$ws.Cells.Item(86, 11).Value = 1898.3265
$ws.Cells.Item(86, 13).Value = -775.3264999999999
$ws.Cells.Item(89, 8).Value = 1896.43
$ws.Cells.Item(89, 9).Value = 1898.3265
$ws.Cells.Item(89, 11).Value = 9491.6325
$ws.Cells.Item(89, 13).Value = -3875.6325
$ws.Cells.Item(134, 8).Value = 15757862
$ws.Cells.Item(134, 9).Value = 20001120
$ws.Cells.Item(134, 10).Value = 3971034.8
$ws.Cells.Item(134, 11).Value = 60003360
$ws.Cells.Item(134, 12).Value = 11913104.4
$ws.Cells.Item(134, 13).Value = -60000825
$ws.Cells.Item(134, 14).Value = -11918174.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3724172.2
$ws.Cells.Item(31, 9).Value = 1895764.8
$ws.Cells.Item(31, 10).Value = 10428333
$ws.Cells.Item(31, 11).Value = 1895764.8
$ws.Cells.Item(31, 12).Value = 10428333
$ws.Cells.Item(31, 13).Value = -1895469.8
$ws.Cells.Item(31, 14).Value = -10428923
$ws.Cells.Item(34, 8).Value = 3724172.2
$ws.Cells.Item(34, 9).Value = 1895764.8
$ws.Cells.Item(34, 10).Value = 10428333
$ws.Cells.Item(34, 11).Value = 1895764.8
$ws.Cells.Item(34, 12).Value = 10428333
$ws.Cells.Item(34, 13).Value = -1895562.8
$ws.Cells.Item(34, 14).Value = -10428737
$ws.Cells.Item(58, 8).Value = 4275950
$ws.Cells.Item(58, 9).Value = 2385492.8
$ws.Cells.Item(58, 10).Value = 11365164
$ws.Cells.Item(58, 11).Value = 2385492.8
$ws.Cells.Item(58, 12).Value = 11365164
$ws.Cells.Item(58, 13).Value = -2385289.8
$ws.Cells.Item(58, 14).Value = -11365570
$ws.Cells.Item(132, 8).Value = 1251916.9
$ws.Cells.Item(132, 9).Value = 1516444.2
$ws.Cells.Item(132, 10).Value = 4858.857
$ws.Cells.Item(132, 11).Value = 4549332.6
$ws.Cells.Item(132, 12).Value = 14576.571
$ws.Cells.Item(132, 13).Value = -4546802.6
$ws.Cells.Item(132, 14).Value = -19636.571
$ws.Cells.Item(134, 8).Value = 1057106
$ws.Cells.Item(134, 9).Value = 4466.4136
$ws.Cells.Item(134, 10).Value = 4448945
$ws.Cells.Item(134, 11).Value = 13399.2408
$ws.Cells.Item(134, 12).Value = 13346835
$ws.Cells.Item(134, 13).Value = -10864.2408
$ws.Cells.Item(134, 14).Value = -13351905
$ws.Cells.Item(136, 8).Value = 4275950
$ws.Cells.Item(136, 9).Value = 2385492.8
$ws.Cells.Item(136, 10).Value = 11365164
$ws.Cells.Item(136, 11).Value = 7156478.399999999
$ws.Cells.Item(136, 12).Value = 34095492
$ws.Cells.Item(136, 13).Value = -7153928.399999999
$ws.Cells.Item(136, 14).Value = -34100592

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1085.125
$ws.Cells.Item(131, 10).Value = 1277.6666
$ws.Cells.Item(131, 12).Value = 3832.9998
$ws.Cells.Item(131, 14).Value = -13912.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 28212752
$ws.Cells.Item(132, 9).Value = 35375064
$ws.Cells.Item(132, 10).Value = 18185518
$ws.Cells.Item(132, 11).Value = 106125192
$ws.Cells.Item(132, 12).Value = 54556554
$ws.Cells.Item(132, 13).Value = -106122662
$ws.Cells.Item(132, 14).Value = -54561614

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2045.125
$ws.Cells.Item(40, 9).Value = 1200
$ws.Cells.Item(40, 11).Value = 1200
$ws.Cells.Item(40, 13).Value = -1064
$ws.Cells.Item(93, 8).Value = 14496.55
$ws.Cells.Item(93, 9).Value = 5908.8335
$ws.Cells.Item(93, 10).Value = 18177
$ws.Cells.Item(93, 11).Value = 5908.8335
$ws.Cells.Item(93, 12).Value = 18177
$ws.Cells.Item(93, 13).Value = -4660.8335
$ws.Cells.Item(93, 14).Value = -20673
$ws.Cells.Item(132, 8).Value = 3179202.8
$ws.Cells.Item(132, 9).Value = 5130097
$ws.Cells.Item(132, 11).Value = 15390291
$ws.Cells.Item(132, 13).Value = -15387761
$ws.Cells.Item(136, 8).Value = 2675729
$ws.Cells.Item(136, 9).Value = 3461768.5
$ws.Cells.Item(136, 10).Value = 3194.8
$ws.Cells.Item(136, 11).Value = 10385305.5
$ws.Cells.Item(136, 12).Value = 9584.400000000001
$ws.Cells.Item(136, 13).Value = -10382755.5
$ws.Cells.Item(136, 14).Value = -14684.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 29960
$ws.Cells.Item(16, 10).Value = 29960
$ws.Cells.Item(16, 12).Value = 29960
$ws.Cells.Item(16, 14).Value = -30544
$ws.Cells.Item(122, 8).Value = 1292.5217
$ws.Cells.Item(122, 9).Value = 1106.9412
$ws.Cells.Item(122, 10).Value = 1818.3334
$ws.Cells.Item(122, 11).Value = 3320.8236
$ws.Cells.Item(122, 12).Value = 5455.0002
$ws.Cells.Item(122, 13).Value = -870.8235999999997
$ws.Cells.Item(122, 14).Value = -10355.0002
$ws.Cells.Item(132, 8).Value = 620824.2
$ws.Cells.Item(132, 9).Value = 2099.6135
$ws.Cells.Item(132, 10).Value = 6065600.5
$ws.Cells.Item(132, 11).Value = 6298.8405
$ws.Cells.Item(132, 12).Value = 18196801.5
$ws.Cells.Item(132, 13).Value = -3768.8405
$ws.Cells.Item(132, 14).Value = -18201861.5
$ws.Cells.Item(136, 8).Value = 8756.966
$ws.Cells.Item(136, 9).Value = 6825.1577
$ws.Cells.Item(136, 10).Value = 12427.4
$ws.Cells.Item(136, 11).Value = 20475.4731
$ws.Cells.Item(136, 12).Value = 37282.2
$ws.Cells.Item(136, 13).Value = -17925.4731
$ws.Cells.Item(136, 14).Value = -42382.2

Write-Host "edit.ps1 completed: applied 172 cell updates across 8 sheets"
